$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 (id 30, 燕返): update frameCount and motionData JSON
$ws.Range("D34").Value = 11
$ws.Range("E34").Value = '{"0": [{"frameNodeType": "MotionFrameNode_Animation", "frameGuid": "70f0", "guid": "279760", "duration": 12, "animSlot": 0}], "1": [], "2": [], "3": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "7fa0", "guid": "158173", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 0.6, "y": 0.6, "z": 0.6}, "count": -3, "slotIndex": -1, "colorHex": "", "colorHex1": ""}], "4": [{"frameNodeType": "MotionFrameNode_AnimPause", "frameGuid": "cbea", "guid": "279760", "pauseCount": 1}], "5": [], "6": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "5c14", "guid": "297400", "offsetPos": {"x": 0, "y": 0, "z": -80}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#B064EEFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "84ad", "guid": "168834", "offsetPos": {"x": -60, "y": 60, "z": -60}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 2, "y": 2, "z": 2}, "count": -10, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_impulse", "frameGuid": "e419", "offsetPos": {"x": -100, "y": 0, "z": 0}, "forceNum": 3000, "groundFriction": 1, "gravityScale": 1, "gravityTime": 0}, {"frameNodeType": "MotionFrameNode_Sound", "frameGuid": "60e5", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "5142", "type": 0, "offsetLoc": {"x": 100, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 20, "y": 4, "z": 4}, "range": 1, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 2, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 30, "dilationRate": 0, "dilationFrame": 0}], "7": [], "8": [], "9": [{"frameNodeType": "MotionFrameNode_BreakPoint", "frameGuid": "28e3", "breakType": 2}], "10": []}'

# Row 35 (id 31, 巨雷一击): update motionData JSON
$ws.Range("E35").Value = '{"0": [{"frameNodeType": "MotionFrameNode_Animation", "frameGuid": "90fd", "guid": "279614", "duration": 30, "animSlot": 0}, {"frameNodeType": "MotionFrameNode_impulse", "frameGuid": "ecd9", "offsetPos": {"x": 60, "y": 0, "z": 100}, "forceNum": 2000, "groundFriction": 1, "gravityScale": 1, "gravityTime": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "f245", "guid": "163348", "offsetPos": {"x": 0, "y": 0, "z": 30}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": -9, "slotIndex": 16, "colorHex": "", "colorHex1": "Color|#F4FF18FF"}], "3": [], "5": [{"frameNodeType": "MotionFrameNode_AnimPause", "frameGuid": "cd41", "guid": "279614", "pauseCount": 2}], "7": [], "8": [{"frameNodeType": "MotionFrameNode_impulse", "frameGuid": "8745", "offsetPos": {"x": 0, "y": 0, "z": 1000}, "forceNum": -6000, "groundFriction": 5, "gravityScale": 1, "gravityTime": 0}, {"frameNodeType": "MotionFrameNode_Sound", "frameGuid": "a47a", "guid": "206115", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000}], "9": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "1bba", "guid": "295657", "offsetPos": {"x": 0, "y": 0, "z": -90}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 2, "y": 2, "z": 2}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "8f3f", "guid": "297401", "offsetPos": {"x": 0, "y": 0, "z": -90}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1.5, "y": 1.5, "z": 1.5}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "83cd", "guid": "297401", "offsetPos": {"x": 0, "y": 0, "z": -90}, "offsetRotation": {"x": 0, "y": 0, "z": 90}, "offsetScale": {"x": 1.5, "y": 1.5, "z": 1.5}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "c71f", "guid": "297401", "offsetPos": {"x": 0, "y": 0, "z": -90}, "offsetRotation": {"x": 0, "y": 0, "z": 180}, "offsetScale": {"x": 1.5, "y": 1.5, "z": 1.5}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "950a", "guid": "297401", "offsetPos": {"x": 0, "y": 0, "z": -90}, "offsetRotation": {"x": 0, "y": 0, "z": 270}, "offsetScale": {"x": 1.5, "y": 1.5, "z": 1.5}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "b4ab", "guid": "180505", "offsetPos": {"x": 0, "y": 0, "z": -60}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 2, "y": 2, "z": 2}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "74cf", "type": 1, "offsetLoc": {"x": 0, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 20, "checkCount": 3, "checkInterval": 0.2, "keepFrameCount": 2, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 30, "dilationRate": 0, "dilationFrame": 0}], "11": [], "12": [], "14": [{"frameNodeType": "MotionFrameNode_BreakPoint", "frameGuid": "7e7a", "breakType": 1}], "15": [], "16": [], "18": []}'

# Row 36 (id 32, 千鸟): update motionData JSON
$ws.Range("E36").Value = '{"0": [{"frameNodeType": "MotionFrameNode_Animation", "frameGuid": "d536", "guid": "279656", "duration": 16, "animSlot": 0}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "8a78", "type": 0, "offsetLoc": {"x": 0, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 1, "checkCount": 1, "checkInterval": 0, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 46, "dilationRate": 0, "dilationFrame": 0}], "1": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "d90c", "guid": "158173", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 0.5, "y": 0.5, "z": 0.5}, "count": -3, "slotIndex": -1, "colorHex": "", "colorHex1": ""}], "2": [{"frameNodeType": "MotionFrameNode_AnimPause", "frameGuid": "b2a2", "guid": "279656", "pauseCount": 1}], "3": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "236a", "guid": "295658", "offsetPos": {"x": 0, "y": 0, "z": -20}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Shake", "frameGuid": "9f35", "rollAmplitude": 0, "rollFrequency": 0, "rollWaveform": 1, "pitchAmplitude": 0, "pitchFrequency": 0, "pitchWaveform": 1, "yawAmplitude": 0, "yawFrequency": 0, "yawWaveform": 1, "xAmplitude": 20, "xFrequency": 0, "xWaveform": 1, "yAmplitude": 0, "yFrequency": 0, "yWaveform": 1, "zAmplitude": 16, "zFrequency": 3, "zWaveform": 1, "fovAmplitude": 100, "fovFrequency": 10, "fovWaveform": 1, "keepTime": 0.5}], "4": [{"frameNodeType": "MotionFrameNode_Move", "frameGuid": "f36b", "move_offsetPos": {"x": 1000, "y": 0, "z": 0}, "move_isFlash": 0, "move_count": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "c2cb", "guid": "200145", "offsetPos": {"x": -1400, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 0.5, "y": 0.7, "z": 0.7}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Sound", "frameGuid": "9255", "guid": "137559", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "de14", "type": 0, "offsetLoc": {"x": 0, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 15, "y": 4, "z": 4}, "range": 1, "checkCount": 1, "checkInterval": 0.1, "keepFrameCount": 2, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "2794", "guid": "301350", "offsetPos": {"x": -1400, "y": 0, "z": -60}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 2, "y": 2, "z": 2}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}], "6": [{"frameNodeType": "MotionFrameNode_BreakPoint", "frameGuid": "15b4", "breakType": 1}], "7": [], "8": [], "9": [], "10": [], "12": [], "13": [], "14": [], "15": [], "16": [], "17": []}'
